# Replace the Word field `{ m:templateVar.name }` (fldChar begin/instrText/
# fldChar end) with literal run text "{m:templateVar.name}", keeping the
# orange accent-6 color formatting on the "template" and "Var" runs.

$d = $word.ActiveDocument

# Find the paragraph that hosts the field (the query's field code).
$hostParagraph = $null
$paraCount = $d.Paragraphs.Count
for ($i = 1; $i -le $paraCount; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Fields.Count -gt 0) {
        $hostParagraph = $candidate
        break
    }
}

# Capture the paragraph's own identity attributes (rsid*) so the rebuilt
# paragraph keeps them, then remove the field - this deletes the
# begin/instrText/end runs that make up the field code, leaving just the
# paragraph mark behind.
$insertionRange = $hostParagraph.Range
$insertionRange.Collapse(1)  # wdCollapseStart
$field = $hostParagraph.Range.Fields.Item(1)
$field.Delete()

# Rebuild the paragraph content as plain literal text runs, preserving the
# paragraph's original identity attributes plus the colored "template"/"Var"
# runs' character formatting.
$xmlFragment = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidP="00F5495F" w:rsidR="00C52979" w:rsidRDefault="00C52979">
<w:r><w:t>{</w:t></w:r>
<w:r><w:t>m</w:t></w:r>
<w:r><w:t>:</w:t></w:r>
<w:r><w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr><w:t>template</w:t></w:r>
<w:r><w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr><w:t>Var</w:t></w:r>
<w:r><w:t xml:space="preserve">.name}</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

[void]$insertionRange.InsertXML($xmlFragment)
